# Add "Trillion Dollar Coach" to the reading list on the "Completed" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Completed" is the tabSelected / active sheet

$newRow = 27

# Text / string columns
$ws.Cells.Item($newRow, 1).Value = "Trillion Dollar Coach"
$ws.Cells.Item($newRow, 2).Value = "Alan Eagle;Eric Schmidt;Jonathan Rosenberg"
$ws.Cells.Item($newRow, 5).Value = "coaching;management;team building;bill campbell;google;silicon valley"
$ws.Cells.Item($newRow, 6).Value = "Audio"
$ws.Cells.Item($newRow, 7).Value = "5 Hours 40 Mins"

# Date columns: copy the date formatting/style from the row above so the new
# cells reuse the existing date style, then overwrite with the new serial value.
$ws.Cells.Item($newRow - 1, 3).Copy()
$ws.Cells.Item($newRow, 3).PasteSpecial(-4122)
$ws.Cells.Item($newRow, 3).Value = 43879

$ws.Cells.Item($newRow - 1, 4).Copy()
$ws.Cells.Item($newRow, 4).PasteSpecial(-4122)
$ws.Cells.Item($newRow, 4).Value = 43879

$excel.CutCopyMode = 0

# Match the scrolled/selected view state shown after the edit.
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("A28").Select()
